# Auto-generated edit script applying scheduled-runner profit recalculations
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8918.733
$ws.Range("I43").Value = 4439.6
$ws.Range("J43").Value = 11158.3
$ws.Range("K43").Value = 4439.6
$ws.Range("L43").Value = 11158.3
$ws.Range("M43").Value = -4370.6
$ws.Range("N43").Value = -11296.3
$ws.Range("H137").Value = 5885066.5
$ws.Range("I137").Value = 1433.6923
$ws.Range("J137").Value = 25006872
$ws.Range("K137").Value = 4301.0769
$ws.Range("L137").Value = 75020616
$ws.Range("M137").Value = -1751.0769
$ws.Range("N137").Value = -75025716
$ws.Range("H138").Value = 5557244.5
$ws.Range("I138").Value = 1572.0555
$ws.Range("J138").Value = 27779936
$ws.Range("K138").Value = 4716.166499999999
$ws.Range("L138").Value = 83339808
$ws.Range("M138").Value = 423.8335000000006
$ws.Range("N138").Value = -83350088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8476465
$ws.Range("I61").Value = 12196783
$ws.Range("K61").Value = 12196783
$ws.Range("M61").Value = -12196571
$ws.Range("H97").Value = 6423.95
$ws.Range("I97").Value = 8936.154
$ws.Range("K97").Value = 8936.154
$ws.Range("M97").Value = -8440.154
$ws.Range("H122").Value = 3997.55
$ws.Range("I122").Value = 4821.6895
$ws.Range("J122").Value = 1824.8182
$ws.Range("K122").Value = 14465.0685
$ws.Range("L122").Value = 5474.4546
$ws.Range("M122").Value = -12015.0685
$ws.Range("N122").Value = -10374.4546
$ws.Range("H132").Value = 4238748.5
$ws.Range("I132").Value = 5815040.5
$ws.Range("J132").Value = 2464.9375
$ws.Range("K132").Value = 17445121.5
$ws.Range("L132").Value = 7394.8125
$ws.Range("M132").Value = -17442591.5
$ws.Range("N132").Value = -12454.8125
$ws.Range("H135").Value = 46277.5
$ws.Range("J135").Value = 46277.5
$ws.Range("L135").Value = 46277.5
$ws.Range("N135").Value = -56417.5
$ws.Range("H136").Value = 8476465
$ws.Range("I136").Value = 12196783
$ws.Range("K136").Value = 36590349
$ws.Range("M136").Value = -36587799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2806.682
$ws.Range("I105").Value = 918.75
$ws.Range("K105").Value = 918.75
$ws.Range("M105").Value = 828.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851498
$ws.Range("I31").Value = 5410.3125
$ws.Range("J31").Value = 13334490
$ws.Range("K31").Value = 5410.3125
$ws.Range("L31").Value = 13334490
$ws.Range("M31").Value = -5115.3125
$ws.Range("N31").Value = -13335080
$ws.Range("H34").Value = 5851498
$ws.Range("I34").Value = 5410.3125
$ws.Range("J34").Value = 13334490
$ws.Range("K34").Value = 5410.3125
$ws.Range("L34").Value = 13334490
$ws.Range("M34").Value = -5208.3125
$ws.Range("N34").Value = -13334894
$ws.Range("H58").Value = 2406.0605
$ws.Range("J58").Value = 2747.6
$ws.Range("L58").Value = 2747.6
$ws.Range("N58").Value = -3153.6
$ws.Range("H92").Value = 32500
$ws.Range("J92").Value = 32500
$ws.Range("L92").Value = 32500
$ws.Range("N92").Value = -37492
$ws.Range("H134").Value = 351567.4
$ws.Range("I134").Value = 1283.5111
$ws.Range("J134").Value = 1036905.44
$ws.Range("K134").Value = 3850.5333
$ws.Range("L134").Value = 3110716.32
$ws.Range("M134").Value = -1315.5333
$ws.Range("N134").Value = -3115786.32
$ws.Range("H136").Value = 2406.0605
$ws.Range("J136").Value = 2747.6
$ws.Range("L136").Value = 8242.799999999999
$ws.Range("N136").Value = -13342.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 995
$ws.Range("I16").Value = 995
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2985
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2812
$ws.Range("N16").ClearContents()
$ws.Range("H20").Value = 2800
$ws.Range("I20").Value = 1200
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 3600
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -3373
$ws.Range("N20").Value = -9454
$ws.Range("H21").Value = 542.2222
$ws.Range("I21").Value = 482.85715
$ws.Range("J21").Value = 750
$ws.Range("K21").Value = 1448.57145
$ws.Range("L21").Value = 2250
$ws.Range("M21").Value = -1275.57145
$ws.Range("N21").Value = -2596
$ws.Range("H22").Value = 1437.5
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -6338
$ws.Range("H26").Value = 225.61111
$ws.Range("J26").Value = 318
$ws.Range("L26").Value = 954
$ws.Range("N26").Value = -1530
$ws.Range("H27").Value = 1437.5
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -6204
$ws.Range("H29").Value = 402.8
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 453.5
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 1360.5
$ws.Range("M29").Value = -323
$ws.Range("N29").Value = -1914.5
$ws.Range("H68").Value = 746.85266
$ws.Range("I68").Value = 505.9861
$ws.Range("J68").Value = 1500.8695
$ws.Range("K68").Value = 1517.9583
$ws.Range("L68").Value = 4502.6085
$ws.Range("M68").Value = -706.9583
$ws.Range("N68").Value = -6124.6085
$ws.Range("H71").Value = 746.85266
$ws.Range("I71").Value = 505.9861
$ws.Range("J71").Value = 1500.8695
$ws.Range("K71").Value = 4553.8749
$ws.Range("L71").Value = 13507.8255
$ws.Range("M71").Value = -497.8748999999998
$ws.Range("N71").Value = -21619.8255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 38344
$ws.Range("J95").Value = 38344
$ws.Range("L95").Value = 38344
$ws.Range("N95").Value = -43836
$ws.Range("H132").Value = 2977.9302
$ws.Range("I132").Value = 1841.9259
$ws.Range("J132").Value = 4894.9375
$ws.Range("K132").Value = 5525.7777
$ws.Range("L132").Value = 14684.8125
$ws.Range("M132").Value = -2995.7777
$ws.Range("N132").Value = -19744.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1954.4546
$ws.Range("I100").Value = 1682.5
$ws.Range("J100").Value = 2280.8
$ws.Range("K100").Value = 1682.5
$ws.Range("L100").Value = 2280.8
$ws.Range("M100").Value = -1141.5
$ws.Range("N100").Value = -3362.8
$ws.Range("H132").Value = 13898653
$ws.Range("I132").Value = 6790.7617
$ws.Range("J132").Value = 33347260
$ws.Range("K132").Value = 20372.2851
$ws.Range("L132").Value = 100041780
$ws.Range("M132").Value = -17842.2851
$ws.Range("N132").Value = -100046840
